$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reporting dashboard data: populate "Edition" (column D) counts for books
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 3

# Modification of the import method: mark finished rows in Comments (column G)
$ws.Range("G3").Value = "### FINISHED"
$ws.Range("G10").Value = "### FINISHED"
